# The deck ships with two color/font/format themes:
#   ppt/theme/theme1.xml -> used by the (only) slide master   -> currently "Integral" / "Red Violet"
#   ppt/theme/theme2.xml -> used by the notes master          -> currently "Office Theme" / "Office"
#
# The target commit swaps the two themes' contents: theme1.xml ends up holding
# the "Office Theme" colors and theme2.xml ends up holding the "Integral" /
# "Red Violet" colors. The font scheme and format (fill/line/effect) scheme are
# identical between the two themes, so the only substantive difference is the
# 12 scheme colors (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
#
# The exposed PowerPoint object model only lets us edit the color scheme that
# backs the (single) slide master / design, via Master.ColorScheme.Colors(n).
# That collection addresses the clrScheme children in document order
# (1=dk1, 2=lt1, 3=dk2, 4=lt2, 5-10=accent1-6, 11=hlink, 12=folHlink), which we
# use to repaint theme1.xml with the "Office Theme" palette.

function HexToComRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r -bor ($g -shl 8) -bor ($b -shl 16)
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.ColorScheme

# Target palette = the deck's current "Office Theme" (theme2.xml), in
# clrScheme child order.
$officeThemeColors = @(
    "000000", # 1  dk1
    "FFFFFF", # 2  lt1
    "44546A", # 3  dk2
    "E7E6E6", # 4  lt2
    "5B9BD5", # 5  accent1
    "ED7D31", # 6  accent2
    "A5A5A5", # 7  accent3
    "FFC000", # 8  accent4
    "4472C4", # 9  accent5
    "70AD47", # 10 accent6
    "0563C1", # 11 hlink
    "954F72"  # 12 folHlink
)

for ($i = 0; $i -lt $officeThemeColors.Count; $i++) {
    $colorScheme.Colors($i + 1).RGB = HexToComRGB $officeThemeColors[$i]
}
